# Genoom Simpsons Family Tree.docx - apply commit changes:
#  1. Remove the stray _GoBack bookmark around "Logic and Code strategy".
#  2. Split "...use a strategy pattern so..." into three runs, inserting "like ".
#  3. Append a new "For the SQL project" section (with a page break) after the
#     Visio diagram, including a relocated _GoBack bookmark, plus the trailing
#     blank paragraphs seen in the target document.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the bookmarkStart/bookmarkEnd (id 0, name "_GoBack") that sits
#    right before the "Logic and Code strategy" Heading2 run, while keeping
#    the heading text itself intact.
# ---------------------------------------------------------------------------
$headingText = "Logic and Code strategy"
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq $headingText) {
        $start = $p.Range.Start
        $end = $p.Range.End - 1
        $r = $d.Range($start, $end)
        $r.Delete()
        $p2 = $d.Paragraphs($i)
        $p2.Range.InsertBefore($headingText)
        break
    }
}

# ---------------------------------------------------------------------------
# 2) "...we could use a strategy pattern so..." -> insert the word "like "
#    between "strategy " and "pattern", split into three separate runs (as
#    happens when Word records a distinct edit/formatting boundary).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("pattern so we can access the data on SQL Server")
$patternStart = $rng.Start

$insertPoint = $d.Range($patternStart, $patternStart)
$insertPoint.InsertBefore("like ")

$likeRange = $d.Range($patternStart, $patternStart + 5)
$likeRange.Font.Bold = 1
$likeRange.Font.Bold = 0

# ---------------------------------------------------------------------------
# 3) Append the new "For the SQL project" block + surrounding blank
#    paragraphs after the Visio diagram paragraph, right before the sectPr.
# ---------------------------------------------------------------------------
function Add-Para {
    param($afterIndex, $innerXml)
    $cur = $d.Paragraphs($afterIndex)
    $r = $cur.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $newIndex = $afterIndex + 1
    $newP = $d.Paragraphs($newIndex)
    $r2 = $newP.Range
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:p>'
    $r2.InsertXML($xml)
    return $newIndex
}

$idx = $d.Paragraphs.Last.Index

# three blank paragraphs
$idx = Add-Para $idx ""
$idx = Add-Para $idx ""
$idx = Add-Para $idx ""

# paragraph containing only a page break
$idx = Add-Para $idx '<w:r><w:br w:type="page"/></w:r>'

# "For the SQL project: " heading
$idx = Add-Para $idx '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">For the SQL project: </w:t></w:r>'

# "We use Dapper instead of SQL Server,"
$idx = Add-Para $idx '<w:r><w:t>We use Dapper instead of SQL Server,</w:t></w:r>'

# "It's a lightweight OR" + relocated _GoBack bookmark + "M and has better performance than EF."
$apos = [char]0x2019
$bmParaXml = '<w:r><w:t>It' + $apos + 's a lightweight OR</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>M and has better performance than EF.</w:t></w:r>'
$idx = Add-Para $idx $bmParaXml

# nine trailing blank paragraphs
for ($i = 0; $i -lt 9; $i++) {
    $idx = Add-Para $idx ""
}

Write-Output ("Done. Paragraph count now " + $d.Paragraphs.Count)
